# Update "想去人数" (F column) counters to the newly scraped values.
# Sheet "展览" (rows 3-22)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 591
$ws.Range("F4").Value = 887
$ws.Range("F5").Value = 639
$ws.Range("F6").Value = 795
$ws.Range("F7").Value = 364
$ws.Range("F8").Value = 559
$ws.Range("F9").Value = 111
$ws.Range("F10").Value = 1131
$ws.Range("F11").Value = 580
$ws.Range("F12").Value = 343
$ws.Range("F13").Value = 455
$ws.Range("F14").Value = 145
$ws.Range("F15").Value = 296
$ws.Range("F16").Value = 44
$ws.Range("F17").Value = 66
$ws.Range("F18").Value = 528
$ws.Range("F19").Value = 20
$ws.Range("F20").Value = 523
$ws.Range("F21").Value = 13
$ws.Range("F22").Value = 492

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 71
$ws.Range("F5").Value = 96
$ws.Range("F11").Value = 12
$ws.Range("F12").Value = 21

# Sheet "全部类型" (aggregated view of all sheets)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 71
$ws.Range("F7").Value = 591
$ws.Range("F8").Value = 887
$ws.Range("F9").Value = 639
$ws.Range("F10").Value = 795
$ws.Range("F11").Value = 364
$ws.Range("F12").Value = 559
$ws.Range("F13").Value = 111
$ws.Range("F14").Value = 1131
$ws.Range("F15").Value = 580
$ws.Range("F16").Value = 96
$ws.Range("F18").Value = 343
$ws.Range("F19").Value = 455
$ws.Range("F21").Value = 145
$ws.Range("F23").Value = 296
$ws.Range("F24").Value = 44
$ws.Range("F25").Value = 66
$ws.Range("F28").Value = 528
$ws.Range("F29").Value = 12
$ws.Range("F30").Value = 21
$ws.Range("F31").Value = 20
$ws.Range("F32").Value = 523
$ws.Range("F33").Value = 13
$ws.Range("F34").Value = 492
